# Insert a new "budget-type" column (B) into the "Data-wide-value" sheet,
# shifting the existing year columns (2013/2014/2015/2016) one column to
# the right (B:E -> C:F), and fill the new column with the budget-type
# value "budget" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# Shift columns B:E to C:F by inserting a fresh column at B.
$ws.Range("B:B").Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "budget-type"

# Every data row (2 through 112) gets the same budget-type value.
$ws.Range("B2:B112").Value = "budget"
